# Applies the "Updated cryptos list" data refresh (Price / Volume(1h) columns,
# plus the Cronos/ImmutableX row swap at rows 45-46) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "73.786.40"
$ws.Cells.Item(2, 5).Value = "  +7.10%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.622.71"
$ws.Cells.Item(3, 5).Value = "  +7.42%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.03%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'184.36"
$ws.Cells.Item(5, 5).Value = "  +13.50%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'580.74"
$ws.Cells.Item(6, 5).Value = "  +3.58%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.999"
$ws.Cells.Item(7, 5).Value = "  -0.10%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.533"
$ws.Cells.Item(8, 5).Value = "  +4.11%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.198"
$ws.Cells.Item(9, 5).Value = "  +17.50%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "2.622.67"
$ws.Cells.Item(10, 5).Value = "  +7.47%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +0.09%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  +7.78%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'4.74"
$ws.Cells.Item(13, 5).Value = "  +3.93%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "73.583.85"
$ws.Cells.Item(14, 5).Value = "  +6.97%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'0.0000188"
$ws.Cells.Item(15, 5).Value = "  +6.11%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "3.076.44"
$ws.Cells.Item(16, 5).Value = "  +6.53%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'26.06"
$ws.Cells.Item(17, 5).Value = "  +11.78%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "2.614.01"
$ws.Cells.Item(18, 5).Value = "  +7.01%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'8.96"
$ws.Cells.Item(19, 5).Value = "  +28.96%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'11.83"
$ws.Cells.Item(20, 5).Value = "  +12.00%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'371.61"
$ws.Cells.Item(21, 5).Value = "  +9.55%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'2.25"
$ws.Cells.Item(22, 5).Value = "  +16.55%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'4.07"
$ws.Cells.Item(23, 5).Value = "  +5.72%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  +0.07%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'69.52"
$ws.Cells.Item(25, 5).Value = "  +3.61%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'4.12"
$ws.Cells.Item(26, 5).Value = "  +10.82%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'9.27"
$ws.Cells.Item(27, 5).Value = "  +12.70%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "2.732.00"
$ws.Cells.Item(28, 5).Value = "  +6.31%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'1.00"
$ws.Cells.Item(29, 5).Value = "  +0.17%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "0.0₃0934"
$ws.Cells.Item(30, 5).Value = "  +13.33%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'515.76"
$ws.Cells.Item(31, 5).Value = "  +19.95%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'1.38"
$ws.Cells.Item(32, 5).Value = "  +18.21%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +6.11%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +7.53%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'0.999"
$ws.Cells.Item(35, 5).Value = "  -0.06%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +12.89%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'161.53"
$ws.Cells.Item(37, 5).Value = "  +1.64%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'19.14"
$ws.Cells.Item(38, 5).Value = "  +6.30%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +1.47%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -0.02%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'4.88"
$ws.Cells.Item(41, 5).Value = "  +11.95%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  +9.84%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "'0.325"
$ws.Cells.Item(43, 5).Value = "  +8.56%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "'157.89"
$ws.Cells.Item(44, 5).Value = "  +21.36%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "ImmutableX"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(45, 4).Value = "'1.18"
$ws.Cells.Item(45, 5).Value = "  +9.61%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "Cronos"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(46, 4).Value = "'0.0870"
$ws.Cells.Item(46, 5).Value = "  +21.21%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'2.33"
$ws.Cells.Item(47, 5).Value = "  +13.59%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'38.61"
$ws.Cells.Item(48, 5).Value = "  +3.04%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'3.61"
$ws.Cells.Item(49, 5).Value = "  +7.86%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'0.527"
$ws.Cells.Item(50, 5).Value = "  +9.35%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'20.39"
$ws.Cells.Item(51, 5).Value = "  +20.67%  "
